$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OS")

# Insert a new row for "CPU" right after the header row, pushing every
# existing attribute row down by one (OS NAME, STATUS, ... USERESSCHEDULER).
$ws.Rows.Item(2).Insert()

# Copy formatting (fill/alignment -> style ids 5/6) from the row that is
# now directly below (the old row 2, "OS NAME") so the new row matches the
# rest of the data rows instead of inheriting the header's shaded style.
$ws.Range("A3:D3").Copy()
$ws.Range("A2:D2").PasteSpecial(-4122)

# Populate the new CPU attribute row.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "CPU"
$ws.Range("D2").Value = "Provide a meaningful name for the CPU that you are going to use."
$ws.Range("C2").Value = "Intel_x86_64"

# Renumber the S.No column for every row that shifted down.
for ($i = 3; $i -le 12; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

# Move the active sheet/selection from TASK to OS.
$ws.Activate() | Out-Null
$ws.Range("A10:A12").Select() | Out-Null
